$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (row 1) and sample values (row 2) added to the
# Project Admin test data sheet, continuing on from column AU.
$ws.Range("AV1").Value = "LSF_LessThan"
$ws.Range("AV2").Value = 1
$ws.Range("AW1").Value = "LSF_GThan"
$ws.Range("AW2").Value = 1
$ws.Range("AX1").Value = "LSF_Per"
$ws.Range("AX2").Value = 1
$ws.Range("AY1").Value = "LSF_Dollar"
$ws.Range("AY2").Value = 1
$ws.Range("AZ1").Value = "CFT_LessThan"
$ws.Range("AZ2").Value = 1
$ws.Range("BA1").Value = "CFT_GThan"
$ws.Range("BA2").Value = 1
$ws.Range("BB1").Value = "CFT_Per"
$ws.Range("BB2").Value = 1
$ws.Range("BC1").Value = "CFT_Dollar"
$ws.Range("BC2").Value = 1
$ws.Range("BD1").Value = "VH_Per"
$ws.Range("BD2").Value = 1
$ws.Range("BE1").Value = "VH_Dol"
$ws.Range("BE2").Value = 1
$ws.Range("BF1").Value = "CF_Per"
$ws.Range("BF2").Value = 1
$ws.Range("BG1").Value = "CF_Dol"
$ws.Range("BG2").Value = 1
$ws.Range("BH1").Value = "VCF_Per"
$ws.Range("BH2").Value = 1
$ws.Range("BI1").Value = "VCF_Dol"
$ws.Range("BI2").Value = 1
$ws.Range("BJ1").Value = "Other_fee_fld1"
$ws.Range("BJ2").Value = 1
$ws.Range("BK1").Value = "Other_fee_fld2"
$ws.Range("BK2").Value = 1
$ws.Range("BL1").Value = "Other_fee_fld3"
$ws.Range("BL2").Value = 1
$ws.Range("BM1").Value = "SCH_LessThan"
$ws.Range("BM2").Value = 1
$ws.Range("BN1").Value = "SCH_GThan"
$ws.Range("BN2").Value = 1
$ws.Range("BO1").Value = "SCH_Per"
$ws.Range("BO2").Value = 1
$ws.Range("BP1").Value = "SCH_Dollar"
$ws.Range("BP2").Value = 1
$ws.Range("BQ1").Value = "Mileage_term_fld1"
$ws.Range("BQ2").Value = 1
$ws.Range("BR1").Value = "Mileage_term_fld2"
$ws.Range("BR2").Value = 1
$ws.Range("BS1").Value = "Mileage_term_fld3"
$ws.Range("BS2").Value = 1
$ws.Range("BT1").Value = "Circuit_fld1"
$ws.Range("BT2").Value = 1
$ws.Range("BU1").Value = "THT_LessThan"
$ws.Range("BU2").Value = 2
$ws.Range("BV1").Value = "THT_GThan"
$ws.Range("BV2").Value = 1
$ws.Range("BW1").Value = "THT_Per"
$ws.Range("BW2").Value = 1
$ws.Range("BX1").Value = "THT_Dollar"
$ws.Range("BX2").Value = 1

# Match the updated view state left in the workbook: scrolled right to the
# new columns, with BP5 as the active selection.
$ws.Range("BP5").Select()
